# Apply the FHIR StructureDefinition metadata refresh (5.0.0 -> 6.0.0 publish):
#  - Metadata sheet: bump Version + Date, add Publisher name, insert a
#    Jurisdiction row, and drop the now-removed duplicate Contact row.
#  - Elements sheet: populate the root Extension's Short/Definition with the
#    real title & description instead of the generic placeholders.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata"
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refreshed publish timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Rows 9-20 get rewritten in place (this both fills in the Publisher name and
# shifts everything below it down by one logical row to make room for the
# new Jurisdiction row), and the old trailing row 21 is removed afterwards.
$meta.Range("A9").Value  = "Publisher"
$meta.Range("B9").Value  = "Alvearie Team"

$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

$meta.Range("A11").Value = "Description"
$meta.Range("B11").Value = "Elimination or waiting period for the short-term disability (STD) benefit (for example, 90, 180 or 365 days). This is the amount of time between the first absent date and the coverage begin date."

$meta.Range("A12").Value = "Purpose"
$meta.Range("B12").Value = ""

$meta.Range("A13").Value = "Copyright"
$meta.Range("B13").Value = ""

$meta.Range("A14").Value = "FHIR Version"
$meta.Range("B14").Value = "4.0.1"

$meta.Range("A15").Value = "Kind"
$meta.Range("B15").Value = "complex-type"

$meta.Range("A16").Value = "Type"
$meta.Range("B16").Value = "Extension"

$meta.Range("A17").Value = "Base Definition"
$meta.Range("B17").Value = "http://hl7.org/fhir/StructureDefinition/Extension"

$meta.Range("A18").Value = "Abstract"
# Leading apostrophe forces this to stay plain text "false" rather than
# being auto-coerced into the boolean FALSE value.
$meta.Range("B18").Value = "'false"

$meta.Range("A19").Value = "Derivation"
$meta.Range("B19").Value = "constraint"

$meta.Range("A20").Value = "Context"
$meta.Range("B20").Value = "element:Element"

# The old row 21 ("Context" / "element:Element") is now a duplicate of row 20
# above, so drop it - this also shrinks the sheet dimension back to B20.
$meta.Rows.Item(21).Delete()

# ---------------------------------------------------------------------------
# Sheet "Elements"
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: replace the generic placeholders with the real
# title/description for this profile.
$elements.Range("K2").Value = "Shortterm Care Waiting Period"
$elements.Range("L2").Value = "Elimination or waiting period for the short-term disability (STD) benefit (for example, 90, 180 or 365 days). This is the amount of time between the first absent date and the coverage begin date."
